$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet lists pharmacy products sorted alphabetically by name (column C).
# Two new products were added, which requires splicing two new rows into the
# existing alphabetical list at the correct position, with every row below
# shifting down by one (the running serial number in column A simply stays a
# contiguous 1..45 count, i.e. always equal to row-6):
#   1) "PHYTO K 10 MG 50 F.C.TAB."   -> belongs between PARAMOL (row 36)
#                                       and PONOFORTE (row 37)
#   2) "VONASPIRE 20MG 14 F.C. TAB"  -> belongs between VITATRON (row 41)
#                                       and ZURCAL (row 42)
# ---------------------------------------------------------------------------

# --- Insert "PHYTO K 10 MG 50 F.C.TAB." as the new row 37 -----------------
# Insert a blank row at 37; this pushes PONOFORTE (formerly row 37) down to
# row 38. Copy PONOFORTE's formatting/merges back up into the new blank row
# 37 so both rows share identical layout, then overwrite each row's data.
$ws.Rows("37:37").Insert()
$ws.Range("A38:Q38").Copy($ws.Range("A37:Q37"))

$ws.Range("C37").Value = "PHYTO K 10 MG 50 F.C.TAB."
$ws.Range("H37").Value = "'0:1"
$ws.Range("L37").Value = "'1"
$ws.Range("N37").Value = "'72.50"
$ws.Range("P37").Value = "'14.5000"
$ws.Range("Q37").Value = "'0:1"

# --- Insert "VONASPIRE 20MG 14 F.C. TAB" as the new row 43 -----------------
# After the previous insert, VITATRON is row 42 and ZURCAL is row 43. Insert
# a blank row at 43; this pushes ZURCAL down to row 44. Copy ZURCAL's
# formatting/merges back up into the new blank row 43, then overwrite each
# row's data.
$ws.Rows("43:43").Insert()
$ws.Range("A44:Q44").Copy($ws.Range("A43:Q43"))

$ws.Range("C43").Value = "VONASPIRE 20MG 14 F.C. TAB"
$ws.Range("H43").Value = "1:0"
$ws.Range("L43").Value = "1"
$ws.Range("N43").Value = "104.00"
$ws.Range("P43").Value = "52.0000"
$ws.Range("Q43").Value = "0:1"

# --- Re-number the serial-number column so it stays a contiguous 1..45 ----
for ($r = 37; $r -le 51; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# --- Update the grand total (sum of the "sell price" column) --------------
# Total grows by the two new sell prices: 14.5000 + 52.0000 = 66.50
$ws.Range("P52").Value = 2433.7800000000002

# --- Update the generated timestamp in the footer --------------------------
$ws.Range("A53").Value = "Monday, 11 August, 2025 1:42 PM"
